$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'60.701.72"
$ws.Range("E2").Value = "  -0.60%  "

# Row 4
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
$ws.Range("D5").Value = "'514.00"
$ws.Range("E5").Value = "  +0.83%  "

# Row 6
$ws.Range("D6").Value = "'154.34"
$ws.Range("E6").Value = "  -2.86%  "

# Row 7
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.19%  "

# Row 8
$ws.Range("D8").Value = "'0.587"
$ws.Range("E8").Value = "  -2.39%  "

# Row 9
$ws.Range("D9").Value = "'2.629.83"
$ws.Range("E9").Value = "  -1.59%  "

# Row 10
$ws.Range("D10").Value = "'6.79"
$ws.Range("E10").Value = "  +4.58%  "

# Row 11
$ws.Range("D11").Value = "'0.105"
$ws.Range("E11").Value = "  -0.89%  "

# Row 12
$ws.Range("E12").Value = "  -0.10%  "

# Row 13
$ws.Range("E13").Value = "  +1.44%  "

# Row 14
$ws.Range("D14").Value = "'3.074.37"
$ws.Range("E14").Value = "  -1.94%  "

# Row 15
$ws.Range("D15").Value = "'60.725.48"
$ws.Range("E15").Value = "  -0.54%  "

# Row 16
$ws.Range("D16").Value = "'21.64"
$ws.Range("E16").Value = "  -1.12%  "

# Row 17
$ws.Range("E17").Value = "  -0.63%  "

# Row 18
$ws.Range("D18").Value = "'2.631.50"
$ws.Range("E18").Value = "  -1.43%  "

# Row 19
$ws.Range("D19").Value = "'4.75"
$ws.Range("E19").Value = "  -1.39%  "

# Row 20
$ws.Range("D20").Value = "'355.14"
$ws.Range("E20").Value = "  +1.89%  "

# Row 21
$ws.Range("D21").Value = "'10.58"
$ws.Range("E21").Value = "  -0.05%  "

# Row 22
$ws.Range("E22").Value = "  -1.01%  "

# Row 23
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  -0.05%  "

# Row 24
$ws.Range("D24").Value = "'61.03"
$ws.Range("E24").Value = "  +0.37%  "

# Row 25
$ws.Range("E25").Value = "  -0.78%  "

# Row 26
$ws.Range("E26").Value = "  -1.60%  "

# Row 27
$ws.Range("D27").Value = "'0.995"
$ws.Range("E27").Value = "  -0.15%  "

# Row 28
$ws.Range("D28").Value = "'0.0₃0843"
$ws.Range("E28").Value = "  -3.55%  "

# Row 29
$ws.Range("E29").Value = "  -4.01%  "

# Row 30
$ws.Range("D30").Value = "'1.00"

# Row 31
$ws.Range("D31").Value = "'19.43"
$ws.Range("E31").Value = "  -0.86%  "

# Row 32
$ws.Range("D32").Value = "'151.93"
$ws.Range("E32").Value = "  -2.95%  "

# Row 33
$ws.Range("D33").Value = "'1.58"
$ws.Range("E33").Value = "  -0.34%  "

# Row 34
$ws.Range("D34").Value = "'5.84"
$ws.Range("E34").Value = "  -0.38%  "

# Row 35
$ws.Range("E35").Value = "  -2.53%  "

# Row 36
$ws.Range("E36").Value = "  -3.05%  "

# Row 37
$ws.Range("D37").Value = "'0.865"
$ws.Range("E37").Value = "  +2.42%  "

# Row 38
$ws.Range("E38").Value = "  -1.33%  "

# Row 39
$ws.Range("E39").Value = "  +2.54%  "

# Row 40
$ws.Range("E40").Value = "  -2.38%  "

# Row 41
$ws.Range("E41").Value = "  -0.83%  "

# Row 42
$ws.Range("D42").Value = "'292.84"
$ws.Range("E42").Value = "  -7.67%  "

# Row 43
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").Value = "'0.101"
$ws.Range("E43").Value = "  +0.67%  "

# Row 44
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").Value = "'0.626"
$ws.Range("E44").Value = "  -2.05%  "

# Row 45
$ws.Range("D45").Value = "'0.996"
$ws.Range("E45").Value = "  -0.43%  "

# Row 46
$ws.Range("D46").Value = "'0.0554"
$ws.Range("E46").Value = "  -4.65%  "

# Row 47
$ws.Range("D47").Value = "'19.83"
$ws.Range("E47").Value = "  +0.26%  "

# Row 48
$ws.Range("D48").Value = "'4.92"
$ws.Range("E48").Value = "  +0.10%  "

# Row 49
$ws.Range("E49").Value = "  -1.08%  "

# Row 50
$ws.Range("E50").Value = "  -0.06%  "

# Row 51
$ws.Range("D51").Value = "'2.004.00"
